# Updated symbol list on Tue Dec 13 23:35:40 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. All cells in this sheet store
# plain text (inline strings) even when the text looks numeric (e.g. "271.80"),
# so we must force text entry (apostrophe prefix) and then clear the resulting
# quote-prefix style so the cell style index is unchanged from the original.
$updates = [ordered]@{
    'D2' = '272.26'
    'D3' = '22.84'
    'D4' = '6.480'
    'D5' = '0.06212'
    'D7' = '6.658'
    'D8' = '1.381'
    'D9' = '0.8294'
    'D10' = '0.01380'
    'D11' = '0.1599'
    'D12' = '0.08276'
    'D14' = '0.03187'
    'D16' = '3.836'
    'D17' = '0.001636'
    'D18' = '0.04732'
    'D19' = '0.006277'
    'D24' = '2.399'
    'D25' = '0.3345'
    'D27' = '0.0002702'
    'D40' = '0.04696'
    'D41' = '0.007062'
    'B42' = 'CEJI'
    'C42' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D42' = '0.003797'
    'E42' = '41CEJICEJI'
    'B43' = 'BKEXToken'
    'C43' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D43' = '0.1162'
    'E43' = '42BKEXTokenBKK'
    'D44' = '0.01203'
    'D45' = '0.00006261'
    'E46' = '45ACDXExchangeACXT'
    'D47' = '0.00000000749'
    'D48' = '0.9194'
    'D49' = '0.002134'
    'E49' = '48BOLOBOLOBestin24h'
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}
